$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.131.77"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.51%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.357.78"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.12%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.683"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.59%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "240.18"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.79%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.93"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +2.18%  "

$ws.Range("E8").Value = "  -0.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.603"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +6.65%  "

$ws.Range("E10").Value = "  +2.55%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.20"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.09%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "32.46"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +15.35%  "

$ws.Range("E13").Value = "  +7.40%  "

$ws.Range("E14").Value = "  +0.98%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.705.73"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.09%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "16.58"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.19%  "

$ws.Range("E17").Value = "  +2.21%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.355.44"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.33%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "44.034.01"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.35%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000103"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.40%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.77"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +5.64%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "77.16"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "257.61"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.32%  "

$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("E25").Value = "  +17.71%  "

$ws.Range("E26").Value = "  -1.29%  "

$ws.Range("E27").Value = "  +0.55%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.77"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.84%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.92"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.77%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.25"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.02%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.26"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.41%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.139"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +4.13%  "

$ws.Range("E33").Value = "  -2.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0764"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +6.08%  "

$ws.Range("E35").Value = "  +2.13%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.38"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.47%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.74"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.98%  "

$ws.Range("E38").Value = "  -2.82%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.37"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.77%  "

$ws.Range("E40").Value = "  +4.30%  "

$ws.Range("E41").Value = "  +11.82%  "

$ws.Range("E42").Value = "  +11.25%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.36"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.65%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.05"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.81%  "

$ws.Range("E45").Value = "  -0.16%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.79"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +7.81%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.52"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +8.85%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.26"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.83%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "101.31"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.36%  "

$ws.Range("E50").Value = "  +1.24%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "57.14"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +9.52%  "
